$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 16.53990612609962
$ws.Range("C2").Value = 5.341090415766084
$ws.Range("E2").Value = 10.28162005241925
$ws.Range("F2").Value = 51.89191793596606
$ws.Range("G2").Value = 3.770587296363396
$ws.Range("I2").Value = 38.73734811103076
$ws.Range("J2").Value = 10.86875401971438
$ws.Range("K2").Value = 14.9147871968598
$ws.Range("L2").Value = 11.20653685119239
$ws.Range("B3").Value = 16.4763354816007
$ws.Range("C3").Value = 5.163190754115814
$ws.Range("E3").Value = 10.31819399675949
$ws.Range("F3").Value = 51.67284419985511
$ws.Range("G3").Value = 3.773697749285053
$ws.Range("I3").Value = 38.61259272432469
$ws.Range("J3").Value = 10.86980068806209
$ws.Range("K3").Value = 14.86473886782259
$ws.Range("L3").Value = 11.23240924027923
$ws.Range("B4").Value = 16.442286713478
$ws.Range("C4").Value = 5.05221400526933
$ws.Range("E4").Value = 10.34233555126476
$ws.Range("F4").Value = 51.54605337122342
$ws.Range("G4").Value = 3.775706932531073
$ws.Range("I4").Value = 38.53995367796106
$ws.Range("J4").Value = 10.87059903977818
$ws.Range("K4").Value = 14.83800408553637
$ws.Range("L4").Value = 11.25043102858891
$ws.Range("B5").Value = 16.42967805327109
$ws.Range("C5").Value = 5.006625810445548
$ws.Range("E5").Value = 10.35259812653061
$ws.Range("F5").Value = 51.49635377266429
$ws.Range("G5").Value = 3.776550763071794
$ws.Range("I5").Value = 38.51135881749826
$ws.Range("J5").Value = 10.87096345972751
$ws.Range("K5").Value = 14.8281230616153
$ws.Range("L5").Value = 11.25831233879044
$ws.Range("B6").Value = 16.42766124034663
$ws.Range("C6").Value = 4.999036216465401
$ws.Range("E6").Value = 10.35432789896748
$ws.Range("F6").Value = 51.48822084109865
$ws.Range("G6").Value = 3.776692397397081
$ws.Range("I6").Value = 38.50667170561239
$ws.Range("J6").Value = 10.87102633008574
$ws.Range("K6").Value = 14.8265437839905
$ws.Range("L6").Value = 11.25965347855335
$ws.Range("B7").Value = 16.44211152454353
$ws.Range("C7").Value = 5.051600559030519
$ws.Range("E7").Value = 10.34247223504306
$ws.Range("F7").Value = 51.5453751000183
$ws.Range("G7").Value = 3.775718211091598
$ws.Range("I7").Value = 38.53956395140484
$ws.Range("J7").Value = 10.87060379631308
$ws.Range("K7").Value = 14.83786671154385
$ws.Range("L7").Value = 11.25053514308071
$ws.Range("B8").Value = 16.51696105344349
$ws.Range("C8").Value = 5.28015597465274
$ws.Range("E8").Value = 10.29388161903931
$ws.Range("F8").Value = 51.81479570603371
$ws.Range("G8").Value = 3.771639213604179
$ws.Range("I8").Value = 38.6935165898944
$ws.Range("J8").Value = 10.86908255780329
$ws.Range("K8").Value = 14.89670713747895
$ws.Range("L8").Value = 11.21501443079414
$ws.Range("B9").Value = 16.70262060525447
$ws.Range("C9").Value = 5.711512550764843
$ws.Range("E9").Value = 10.2119214133123
$ws.Range("F9").Value = 52.40310441187912
$ws.Range("G9").Value = 3.764424566626125
$ws.Range("I9").Value = 39.02637273101814
$ws.Range("J9").Value = 10.86733781623816
$ws.Range("K9").Value = 15.0433389815875
$ws.Range("L9").Value = 11.16229813690918
$ws.Range("B10").Value = 16.86172218517065
$ws.Range("C10").Value = 6.014578174507192
$ws.Range("E10").Value = 10.1597692136934
$ws.Range("F10").Value = 52.87003371300219
$ws.Range("G10").Value = 3.75959637919079
$ws.Range("I10").Value = 39.28910881107585
$ws.Range("J10").Value = 10.86681529341121
$ws.Range("K10").Value = 15.16943059693273
$ws.Range("L10").Value = 11.13387992979512
$ws.Range("B11").Value = 16.93878358901563
$ws.Range("C11").Value = 6.148801867186973
$ws.Range("E11").Value = 10.13778219371191
$ws.Range("F11").Value = 53.08954608990316
$ws.Range("G11").Value = 3.757501278598014
$ws.Range("I11").Value = 39.41241918335297
$ws.Range("J11").Value = 10.86674339126976
$ws.Range("K11").Value = 15.23060974972751
$ws.Range("L11").Value = 11.12318655659017
$ws.Range("B12").Value = 16.96861561989735
$ws.Range("C12").Value = 6.199059241766834
$ws.Range("E12").Value = 10.12970510695006
$ws.Range("F12").Value = 53.17365083110202
$ws.Range("G12").Value = 3.756722387773222
$ws.Range("I12").Value = 39.4596431292692
$ws.Range("J12").Value = 10.86674007275074
$ws.Range("K12").Value = 15.25430997865764
$ws.Range("L12").Value = 11.11945802319707
$ws.Range("B13").Value = 16.96216220721701
$ws.Range("C13").Value = 6.188261461910239
$ws.Range("E13").Value = 10.13143359555045
$ws.Range("F13").Value = 53.15549445373472
$ws.Range("G13").Value = 3.756889493248904
$ws.Range("I13").Value = 39.44944934538994
$ws.Range("J13").Value = 10.86673972301815
$ws.Range("K13").Value = 15.24918227278868
$ws.Range("L13").Value = 11.12024676946974
$ws.Range("B14").Value = 16.94122498737988
$ws.Range("C14").Value = 6.152948247251374
$ws.Range("E14").Value = 10.13711270331164
$ws.Range("F14").Value = 53.09644603656201
$ws.Range("G14").Value = 3.757436909092039
$ws.Range("I14").Value = 39.41629381913383
$ws.Range("J14").Value = 10.86674263869358
$ws.Range("K14").Value = 15.2325489957578
$ws.Range("L14").Value = 11.12287338123022
$ws.Range("B15").Value = 16.9284843470477
$ws.Range("C15").Value = 6.131242334497306
$ws.Range("E15").Value = 10.14062371102439
$ws.Range("F15").Value = 53.06040356592089
$ws.Range("G15").Value = 3.757774100129535
$ws.Range("I15").Value = 39.39605350035437
$ws.Range("J15").Value = 10.86674754030548
$ws.Range("K15").Value = 15.22242954949684
$ws.Range("L15").Value = 11.12452402125209
$ws.Range("B16").Value = 16.85677836772456
$ws.Range("C16").Value = 6.005728954523719
$ws.Range("E16").Value = 10.16124099877181
$ws.Range("F16").Value = 52.85582751542566
$ws.Range("G16").Value = 3.759735329513641
$ws.Range("I16").Value = 39.28112510590448
$ws.Range("J16").Value = 10.86682333448072
$ws.Range("K16").Value = 15.16550788790457
$ws.Range("L16").Value = 11.1346236873012
$ws.Range("B17").Value = 16.81397299373309
$ws.Range("C17").Value = 5.927762915754293
$ws.Range("E17").Value = 10.17433335518477
$ws.Range("F17").Value = 52.73211734274162
$ws.Range("G17").Value = 3.76096435722553
$ws.Range("I17").Value = 39.21157990150931
$ws.Range("J17").Value = 10.8669123450716
$ws.Range("K17").Value = 15.13155531243877
$ws.Range("L17").Value = 11.14139143753676
$ws.Range("B18").Value = 16.78979494967935
$ws.Range("C18").Value = 5.8825785036516
$ws.Range("E18").Value = 10.18202730651109
$ws.Range("F18").Value = 52.66163376751596
$ws.Range("G18").Value = 3.761680798075045
$ws.Range("I18").Value = 39.17193712653027
$ws.Range("J18").Value = 10.86697914318482
$ws.Range("K18").Value = 15.1123872281239
$ws.Range("L18").Value = 11.14549439583704
$ws.Range("B19").Value = 16.78168536067501
$ws.Range("C19").Value = 5.867222875106004
$ws.Range("E19").Value = 10.18466046767075
$ws.Range("F19").Value = 52.63788575831626
$ws.Range("G19").Value = 3.761925013114544
$ws.Range("I19").Value = 39.1585766904884
$ws.Range("J19").Value = 10.86700443735801
$ws.Range("K19").Value = 15.10595964534689
$ws.Range("L19").Value = 11.14691972367419
$ws.Range("B20").Value = 16.81848406849518
$ws.Range("C20").Value = 5.936098117944843
$ws.Range("E20").Value = 10.17292272871115
$ws.Range("F20").Value = 52.74521732496319
$ws.Range("G20").Value = 3.760832538734258
$ws.Range("I20").Value = 39.21894619161704
$ws.Range("J20").Value = 10.86690125455156
$ws.Range("K20").Value = 15.13513242215203
$ws.Range("L20").Value = 11.14064923400405
$ws.Range("B21").Value = 16.94735729156512
$ws.Range("C21").Value = 6.163336414375907
$ws.Range("E21").Value = 10.13543786320214
$ws.Range("F21").Value = 53.11376372467835
$ws.Range("G21").Value = 3.757275727516415
$ws.Range("I21").Value = 39.42601817653781
$ws.Range("J21").Value = 10.86674113286559
$ws.Range("K21").Value = 15.23742026060287
$ws.Range("L21").Value = 11.12209317894653
$ws.Range("B22").Value = 17.03536369679722
$ws.Range("C22").Value = 6.308508150581008
$ws.Range("E22").Value = 10.1123898223145
$ws.Range("F22").Value = 53.36032291437579
$ws.Range("G22").Value = 3.755035497907418
$ws.Range("I22").Value = 39.56442699428281
$ws.Range("J22").Value = 10.86677587579041
$ws.Range("K22").Value = 15.30736990063682
$ws.Range("L22").Value = 11.11183540709148
$ws.Range("B23").Value = 16.98805509715187
$ws.Range("C23").Value = 6.231347030147137
$ws.Range("E23").Value = 10.12455857150417
$ws.Range("F23").Value = 53.22822284076049
$ws.Range("G23").Value = 3.756223459898978
$ws.Range("I23").Value = 39.49027973053135
$ws.Range("J23").Value = 10.86674455651907
$ws.Range("K23").Value = 15.26975853632696
$ws.Range("L23").Value = 11.11713927067181
$ws.Range("B24").Value = 16.81644326629174
$ws.Range("C24").Value = 5.932330892960225
$ws.Range("E24").Value = 10.17355995276043
$ws.Range("F24").Value = 52.73929282925484
$ws.Range("G24").Value = 3.760892103165669
$ws.Range("I24").Value = 39.2156148325026
$ws.Range("J24").Value = 10.86690621990725
$ws.Range("K24").Value = 15.13351411411591
$ws.Range("L24").Value = 11.14098412356594
$ws.Range("B25").Value = 16.64832979123515
$ws.Range("C25").Value = 5.597015168181664
$ws.Range("E25").Value = 10.23267339906114
$ws.Range("F25").Value = 52.23771424689101
$ws.Range("G25").Value = 3.766292949798448
$ws.Range("I25").Value = 38.93307928227087
$ws.Range("J25").Value = 10.86767677634154
$ws.Range("K25").Value = 15.00039239855123
$ws.Range("L25").Value = 11.17474698430932
